# Phase 1: Optimisation UI refinements - remove headings, update labels,
# add comma formatting, fix chart
#
# This script applies four logical changes to the document:
#  1. Extend the "Please explain ... in below" paragraph with additional
#     commentary text (partly in the default black colour, partly in the
#     blue "accent1" colour used elsewhere for author replies), including
#     the "types" grammar-check markers.
#  2. Re-colour a run of bullet points (the "Optimisation UI" list items)
#     in red (EE0000) to flag them for phase-1 attention.
#  3. Move the lastRenderedPageBreak marker from the "Can we split this
#     page..." run up onto the "Output" heading run.

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $matchText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.Contains($matchText)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) Append extra commentary runs after "...in below"
# ---------------------------------------------------------------------
$idxInBelow = Find-ParagraphIndex $d "in below"
$pInBelow = $d.Paragraphs.Item($idxInBelow)
$insertPoint = $pInBelow.Range.End - 1
$insertRange = $d.Range($insertPoint, $insertPoint)

$xmlExplain = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="20A66EC3" w14:textId="1B61DFF3" w:rsidR="00A06A70" w:rsidRPr="0045665E" w:rsidRDefault="00580702" w:rsidP="00EE78EB"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:color w:val="2E74B5" w:themeColor="accent5" w:themeShade="BF"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/></w:rPr><w:t xml:space="preserve">this is to adjust the proportion between accumulation phase and pension phase.  Let me know if you prefer just to have the two different entity </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/></w:rPr><w:t>types</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/></w:rPr><w:t xml:space="preserve"> I thought a slider may help to allocate it better.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertRange.InsertXML($xmlExplain)
# Merge the freshly-inserted paragraph back into the original one by
# removing the paragraph mark that now separates them; the new paragraph
# carries the original paraId/rsid/pPr so the merged paragraph keeps them.
$d.Range($insertPoint, $insertPoint + 1).Delete()

# ---------------------------------------------------------------------
# 2) Colour the "Optimisation UI" bullet list items red (EE0000)
# ---------------------------------------------------------------------
$redBulletMarkers = @(
    "Total Investable Assets",
    "Cashflow Projection Inputs",
    "Projection Period (Years)",
    "Inflation Rate (% p.a.)",
    "Advice Fee (",
    "Make these sections smaller",
    "Income Streams (Today’s Dollars)",
    "Expense Streams (Today’s Dollars)",
    "ongoing or one off"
)

foreach ($marker in $redBulletMarkers) {
    $idx = Find-ParagraphIndex $d $marker
    if ($idx -ge 0) {
        $d.Paragraphs.Item($idx).Range.Font.Color = 238   # BGR for EE0000
    }
}

# ---------------------------------------------------------------------
# 3) Move the lastRenderedPageBreak marker from the "split this page"
#    bullet up onto the "Output" heading run.
# ---------------------------------------------------------------------
$idxOutput = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Output`r") {
        $idxOutput = $i
        break
    }
}
$pOutput = $d.Paragraphs.Item($idxOutput)
$rOutput = $d.Range($pOutput.Range.Start, $pOutput.Range.End - 1)
$xmlOutput = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Output</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rOutput.InsertXML($xmlOutput)

$idxSplit = Find-ParagraphIndex $d "split this page"
$pSplit = $d.Paragraphs.Item($idxSplit)
$rSplit = $d.Range($pSplit.Range.Start, $pSplit.Range.End - 1)
$xmlSplit = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Can we split this page into two with data on the left and charts on the right?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rSplit.InsertXML($xmlSplit)
